{"js": "// Separate pseudocode into different documents:\n// Split the underlined heading \"Country Populations Pseudocode\" into two\n// runs: \"World\" and \" Population Pseudocode\" (both keep the single\n// underline formatting), matching the target OOXML exactly.\n\nconst body = context.document.body;\nconst results = body.search(\"Country Populations Pseudocode\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"Country Populations Pseudocode\" in the document body.');\n}\n\nconst target = results.items[0];\n\n// Use insertOoxml so the replacement produces two distinct <w:r> runs\n// (one for \"World\", one for \" Population Pseudocode\") instead of having\n// the engine merge same-formatted adjacent runs into a single run, which\n// is what a plain insertText() call would do.\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p>' +\n  '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>World</w:t></w:r>' +\n  '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t xml:space=\"preserve\"> Population Pseudocode</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Separate pseudocode into different documents:\n# Split the underlined heading \"Country Populations Pseudocode\" into two\n# runs: \"World\" and \" Population Pseudocode\" (both keep the single\n# underline formatting), matching the target OOXML exactly.\n\n$d = $word.ActiveDocument\n\n# Locate the heading text.\n$range = $d.Content\n$found = $range.Find.Execute(\"Country Populations Pseudocode\")\nif (-not $found) {\n    throw \"Could not find 'Country Populations Pseudocode' in the document.\"\n}\n\n# Capture the host paragraph's own <w:p ...> opening tag (with its\n# w14:paraId / w:rsidR / etc. attributes) so the rebuilt paragraph keeps\n# the exact same paragraph identity instead of getting a bare <w:p>.\n$para = $range.Paragraphs(1)\n$openXml = $para.Range.WordOpenXML\n$pTag = \"<w:p>\"\nif ($openXml -match \"(<w:p[ >][^>]*>)\") {\n    $pTag = $matches[1]\n}\n\n# Build the replacement paragraph: same opening <w:p> tag, but now with two\n# runs - \"World\" and \" Population Pseudocode\" - both carrying the single\n# underline formatting that the original single run had.\n$newParaXml = $pTag +\n    '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t>World</w:t></w:r>' +\n    '<w:r><w:rPr><w:u w:val=\"single\"/></w:rPr><w:t xml:space=\"preserve\"> Population Pseudocode</w:t></w:r>' +\n    '</w:p>'\n\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">' +\n    '<w:body>' + $newParaXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n# Delete the old heading text (collapses the found range to the insertion\n# point where the heading used to live), then insert the two-run XML there.\n# Doing the delete first - instead of inserting then deleting the leftover\n# text - avoids the engine's same-formatting run-merge pass that a\n# post-insert Delete()/Text=\"\" would otherwise trigger.\n$range.Delete()\n$range.InsertXML($xml)\n"}
